$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1814671814671815
$ws.Range("C2").Value = 0.5598455598455598
$ws.Range("J2").Value = 0.007722007722007722
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1081081081081081
$ws.Range("B3").Value = 0.02758620689655172
$ws.Range("C3").Value = 0.02068965517241379
$ws.Range("J3").Value = 0.04137931034482759
$ws.Range("P3").Value = 0.696551724137931
$ws.Range("S3").Value = 0.2137931034482759
$ws.Range("P4").Value = 0.7291666666666666
$ws.Range("S4").Value = 0.2708333333333333
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.07114624505928854
$ws.Range("D6").Value = 0.0158102766798419
$ws.Range("F6").Value = 0.05533596837944664
$ws.Range("J6").Value = 0.2490118577075099
$ws.Range("O6").Value = 0.02371541501976284
$ws.Range("Q6").Value = 0.1343873517786561
$ws.Range("R6").Value = 0.05533596837944664
$ws.Range("S6").Value = 0.3952569169960474
$ws.Range("B7").Value = 0.05660377358490566
$ws.Range("D7").Value = 0.03773584905660377
$ws.Range("F7").Value = 0.09433962264150944
$ws.Range("J7").Value = 0.1037735849056604
$ws.Range("O7").Value = 0.02358490566037736
$ws.Range("Q7").Value = 0.160377358490566
$ws.Range("R7").Value = 0.07075471698113207
$ws.Range("S7").Value = 0.4528301886792453
$ws.Range("B8").Value = 0.06888888888888889
$ws.Range("D8").Value = 0.01777777777777778
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.1311111111111111
$ws.Range("O8").Value = 0.008888888888888889
$ws.Range("Q8").Value = 0.18
$ws.Range("R8").Value = 0.1133333333333333
$ws.Range("S8").Value = 0.4133333333333333
$ws.Range("B9").Value = 0.08673469387755102
$ws.Range("F9").Value = 0.06122448979591837
$ws.Range("J9").Value = 0.09183673469387756
$ws.Range("O9").Value = 0.03061224489795918
$ws.Range("Q9").Value = 0.1785714285714286
$ws.Range("R9").Value = 0.1173469387755102
$ws.Range("S9").Value = 0.4336734693877551
$ws.Range("B10").Value = 0.1063464837049743
$ws.Range("D10").Value = 0.0274442538593482
$ws.Range("E10").Value = 0.0008576329331046312
$ws.Range("F10").Value = 0.07289879931389365
$ws.Range("J10").Value = 0.09691252144082332
$ws.Range("O10").Value = 0.02058319039451115
$ws.Range("Q10").Value = 0.1981132075471698
$ws.Range("R10").Value = 0.08747855917667238
$ws.Range("S10").Value = 0.3893653516295026
$ws.Range("G11").Value = 0.171003717472119
$ws.Range("J11").Value = 0.05947955390334572
$ws.Range("K11").Value = 0.1821561338289963
$ws.Range("L11").Value = 0.5836431226765799
$ws.Range("S11").Value = 0.003717472118959108
$ws.Range("G12").Value = 0.7816091954022989
$ws.Range("J12").Value = 0.132183908045977
$ws.Range("K12").Value = 0.01724137931034483
$ws.Range("L12").Value = 0.05747126436781609
$ws.Range("S12").Value = 0.01149425287356322
$ws.Range("G13").Value = 0.7083333333333334
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.0625
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2
$ws.Range("F15").Value = 0.05181347150259067
$ws.Range("H15").Value = 0.150259067357513
$ws.Range("I15").Value = 0.09844559585492228
$ws.Range("J15").Value = 0.310880829015544
$ws.Range("K15").Value = 0.08808290155440414
$ws.Range("M15").Value = 0.02072538860103627
$ws.Range("N15").Value = 0.005181347150259068
$ws.Range("O15").Value = 0.05699481865284974
$ws.Range("S15").Value = 0.2176165803108808
$ws.Range("F16").Value = 0.005847953216374269
$ws.Range("H16").Value = 0.1695906432748538
$ws.Range("I16").Value = 0.05847953216374269
$ws.Range("J16").Value = 0.4678362573099415
$ws.Range("K16").Value = 0.08771929824561403
$ws.Range("M16").Value = 0.02339181286549707
$ws.Range("O16").Value = 0.04093567251461988
$ws.Range("S16").Value = 0.1461988304093567
$ws.Range("F17").Value = 0.02926829268292683
$ws.Range("H17").Value = 0.2024390243902439
$ws.Range("I17").Value = 0.0975609756097561
$ws.Range("J17").Value = 0.424390243902439
$ws.Range("K17").Value = 0.08048780487804878
$ws.Range("M17").Value = 0.01951219512195122
$ws.Range("O17").Value = 0.03414634146341464
$ws.Range("S17").Value = 0.1121951219512195
$ws.Range("F18").Value = 0.009803921568627451
$ws.Range("H18").Value = 0.1911764705882353
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.4166666666666667
$ws.Range("K18").Value = 0.08333333333333333
$ws.Range("M18").Value = 0.02941176470588235
$ws.Range("O18").Value = 0.08333333333333333
$ws.Range("S18").Value = 0.1029411764705882
$ws.Range("F19").Value = 0.02931854199683043
$ws.Range("H19").Value = 0.2179080824088748
$ws.Range("I19").Value = 0.08637083993660856
$ws.Range("J19").Value = 0.358161648177496
$ws.Range("K19").Value = 0.1014263074484945
$ws.Range("M19").Value = 0.0213946117274168
$ws.Range("N19").Value = 0.003169572107765452
$ws.Range("O19").Value = 0.06339144215530904
$ws.Range("S19").Value = 0.1188589540412044
